$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 22:52"

# Rows 4-4
$block = New-Object 'object[,]' 1,8
$block[0,0] = "Estados Unidos"
$block[0,1] = 306768
$block[0,2] = 29607
$block[0,3] = 14686
$block[0,4] = 283735
$block[0,5] = 8073
$block[0,6] = 943
$block[0,7] = 8347
$ws.Range("A4:H4").Value = $block

# Rows 63-63
$block = New-Object 'object[,]' 1,8
$block[0,0] = "Marruecos"
$block[0,1] = 919
$block[0,2] = 128
$block[0,3] = 66
$block[0,4] = 794
$block[0,5] = 1
$block[0,6] = 11
$block[0,7] = 59
$ws.Range("A63:H63").Value = $block

# Rows 73-78
$block = New-Object 'object[,]' 6,8
$block[0,0] = "Camerun"
$block[0,1] = 555
$block[0,2] = 46
$block[0,3] = 17
$block[0,4] = 529
$block[0,5] = 0
$block[0,6] = 1
$block[0,7] = 9
$block[1,0] = "Tunez"
$block[1,1] = 553
$block[1,2] = 58
$block[1,3] = 5
$block[1,4] = 530
$block[1,5] = 26
$block[1,6] = 0
$block[1,7] = 18
$block[2,0] = "Kazajistan"
$block[2,1] = 531
$block[2,2] = 67
$block[2,3] = 36
$block[2,4] = 490
$block[2,5] = 6
$block[2,6] = 2
$block[2,7] = 5
$block[3,0] = "Azerbaiyan"
$block[3,1] = 521
$block[3,2] = 78
$block[3,3] = 32
$block[3,4] = 484
$block[3,5] = 17
$block[3,6] = 0
$block[3,7] = 5
$block[4,0] = "Libano"
$block[4,1] = 520
$block[4,2] = 12
$block[4,3] = 54
$block[4,4] = 449
$block[4,5] = 27
$block[4,6] = 0
$block[4,7] = 17
$block[5,0] = "Letonia"
$block[5,1] = 509
$block[5,2] = 16
$block[5,3] = 1
$block[5,4] = 507
$block[5,5] = 3
$block[5,6] = 0
$block[5,7] = 1
$ws.Range("A73:H78").Value = $block

# Rows 107-107
$block = New-Object 'object[,]' 1,8
$block[0,0] = "Mauricio"
$block[0,1] = 196
$block[0,2] = 10
$block[0,3] = 7
$block[0,4] = 182
$block[0,5] = 1
$block[0,6] = 0
$block[0,7] = 7
$ws.Range("A107:H107").Value = $block

# Rows 113-121
$block = New-Object 'object[,]' 9,8
$block[0,0] = "Niger"
$block[0,1] = 144
$block[0,2] = 24
$block[0,3] = 0
$block[0,4] = 136
$block[0,5] = 0
$block[0,6] = 3
$block[0,7] = 8
$block[1,0] = "Kirguistan"
$block[1,1] = 144
$block[1,2] = 14
$block[1,3] = 9
$block[1,4] = 134
$block[1,5] = 5
$block[1,6] = 0
$block[1,7] = 1
$block[2,0] = "Martinica"
$block[2,1] = 143
$block[2,2] = 0
$block[2,3] = 27
$block[2,4] = 113
$block[2,5] = 18
$block[2,6] = 0
$block[2,7] = 3
$block[3,0] = "Bolivia"
$block[3,1] = 139
$block[3,2] = 7
$block[3,3] = 1
$block[3,4] = 128
$block[3,5] = 3
$block[3,6] = 1
$block[3,7] = 10
$block[4,0] = "Brunei"
$block[4,1] = 135
$block[4,2] = 1
$block[4,3] = 66
$block[4,4] = 68
$block[4,5] = 3
$block[4,6] = 0
$block[4,7] = 1
$block[5,0] = "Mayotte"
$block[5,1] = 134
$block[5,2] = 6
$block[5,3] = 14
$block[5,4] = 118
$block[5,5] = 3
$block[5,6] = 0
$block[5,7] = 2
$block[6,0] = "Guadalupe"
$block[6,1] = 134
$block[6,2] = 4
$block[6,3] = 24
$block[6,4] = 103
$block[6,5] = 14
$block[6,6] = 0
$block[6,7] = 7
$block[7,0] = "Isla de Man"
$block[7,1] = 126
$block[7,2] = 12
$block[7,3] = 0
$block[7,4] = 125
$block[7,5] = 0
$block[7,6] = 0
$block[7,7] = 1
$block[8,0] = "Kenia"
$block[8,1] = 126
$block[8,2] = 4
$block[8,3] = 4
$block[8,4] = 118
$block[8,5] = 2
$block[8,6] = 0
$block[8,7] = 4
$ws.Range("A113:H121").Value = $block

# Rows 123-128
$block = New-Object 'object[,]' 6,8
$block[0,0] = "Guinea"
$block[0,1] = 111
$block[0,2] = 38
$block[0,3] = 5
$block[0,4] = 106
$block[0,5] = 0
$block[0,6] = 0
$block[0,7] = 0
$block[1,0] = "Trinidad yTobago"
$block[1,1] = 103
$block[1,2] = 5
$block[1,3] = 1
$block[1,4] = 96
$block[1,5] = 0
$block[1,6] = 0
$block[1,7] = 6
$block[2,0] = "Ruanda"
$block[2,1] = 102
$block[2,2] = 13
$block[2,3] = 0
$block[2,4] = 102
$block[2,5] = 0
$block[2,6] = 0
$block[2,7] = 0
$block[3,0] = "Gibraltar"
$block[3,1] = 98
$block[3,2] = 3
$block[3,3] = 52
$block[3,4] = 46
$block[3,5] = 0
$block[3,6] = 0
$block[3,7] = 0
$block[4,0] = "Paraguay"
$block[4,1] = 96
$block[4,2] = 4
$block[4,3] = 12
$block[4,4] = 81
$block[4,5] = 2
$block[4,6] = 0
$block[4,7] = 3
$block[5,0] = "Liechtenstein"
$block[5,1] = 77
$block[5,2] = 2
$block[5,3] = 0
$block[5,4] = 76
$block[5,5] = 0
$block[5,6] = 1
$block[5,7] = 1
$ws.Range("A123:H128").Value = $block

# Rows 144-145
$block = New-Object 'object[,]' 2,8
$block[0,0] = "Puerto Rico"
$block[0,1] = 39
$block[0,2] = 0
$block[0,3] = 1
$block[0,4] = 36
$block[0,5] = 0
$block[0,6] = 0
$block[0,7] = 2
$block[1,0] = "Zambia"
$block[1,1] = 39
$block[1,2] = 0
$block[1,3] = 2
$block[1,4] = 36
$block[1,5] = 0
$block[1,6] = 0
$block[1,7] = 1
$ws.Range("A144:H145").Value = $block

